# Se modifica builder para insertar los datos de la prueba
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the test-data row (row 2) with the new "segundo factor" values.
# Leading apostrophes force text storage so Excel keeps the quotePrefix
# style already applied to these cells (they look like numbers).
$ws.Range("D2").Value = "seguridad01"
$ws.Range("E2").Value = "'1234"
$ws.Range("F2").Value = "'4321"

# Move the active selection to match the saved view state.
$ws.Range("I16").Select()
